$wb = $excel.ActiveWorkbook
$wsPaper = $wb.Worksheets.Item("paper")
$wsForecasted = $wb.Worksheets.Item("forecasted")

# paper sheet (specific consumption per TWh) updates
$wsPaper.Range("B2").Value = 2.3804361
$wsPaper.Range("C2").Value = 6.7645984
$wsPaper.Range("D2").Value = $null
$wsPaper.Range("E2").Value = 0.9183789
$wsPaper.Range("F2").Value = 1.4492306
$wsPaper.Range("G2").Value = 4.396989
$wsPaper.Range("H2").Value = 0
$wsPaper.Range("B3").Value = 0.3329083
$wsPaper.Range("C3").Value = 1.8265926
$wsPaper.Range("D3").Value = $null
$wsPaper.Range("E3").Value = 0.2479828
$wsPaper.Range("F3").Value = 0.3913246
$wsPaper.Range("G3").Value = 1.1872852
$wsPaper.Range("H3").Value = 0
$wsPaper.Range("B4").Value = 2.0434377
$wsPaper.Range("C4").Value = 6.4617311
$wsPaper.Range("D4").Value = $null
$wsPaper.Range("E4").Value = 0.8772608
$wsPaper.Range("F4").Value = 1.384345
$wsPaper.Range("G4").Value = 4.2001252
$wsPaper.Range("H4").Value = 0
$wsPaper.Range("B5").Value = 0.3784892
$wsPaper.Range("D5").Value = $null
$wsPaper.Range("E5").Value = 0.07403709999999999
$wsPaper.Range("F5").Value = 0.1168328
$wsPaper.Range("G5").Value = 0.3544727
$wsPaper.Range("H5").Value = 0
$wsPaper.Range("B6").Value = 23.4268672
$wsPaper.Range("C6").Value = 45.3082845
$wsPaper.Range("D6").Value = $null
$wsPaper.Range("E6").Value = 6.1511664
$wsPaper.Range("F6").Value = 9.706733099999999
$wsPaper.Range("G6").Value = 29.4503849
$wsPaper.Range("H6").Value = 0
$wsPaper.Range("B7").Value = 0.0663697
$wsPaper.Range("C7").Value = 0.0195041
$wsPaper.Range("D7").Value = $null
$wsPaper.Range("E7").Value = 0.0026479
$wsPaper.Range("F7").Value = 0.0041785
$wsPaper.Range("G7").Value = 0.0126777
$wsPaper.Range("H7").Value = 0
$wsPaper.Range("B8").Value = 0.5333476
$wsPaper.Range("C8").Value = 0.6344338
$wsPaper.Range("D8").Value = $null
$wsPaper.Range("E8").Value = 0.08613229999999999
$wsPaper.Range("F8").Value = 0.1359195
$wsPaper.Range("G8").Value = 0.412382
$wsPaper.Range("H8").Value = 0
$wsPaper.Range("B9").Value = 3.8584664
$wsPaper.Range("C9").Value = 10.5076073
$wsPaper.Range("D9").Value = $null
$wsPaper.Range("E9").Value = 1.4265391
$wsPaper.Range("F9").Value = 2.2511234
$wsPaper.Range("G9").Value = 6.8299448
$wsPaper.Range("H9").Value = 0
$wsPaper.Range("B10").Value = 10.5328449
$wsPaper.Range("C10").Value = 18.1152198
$wsPaper.Range("D10").Value = $null
$wsPaper.Range("E10").Value = 2.4593677
$wsPaper.Range("F10").Value = 3.8809592
$wsPaper.Range("G10").Value = 11.7748929
$wsPaper.Range("H10").Value = 0
$wsPaper.Range("B11").Value = 0.2714396
$wsPaper.Range("C11").Value = 0.2435884
$wsPaper.Range("D11").Value = $null
$wsPaper.Range("E11").Value = 0.0330702
$wsPaper.Range("F11").Value = 0.0521858
$wsPaper.Range("G11").Value = 0.1583325
$wsPaper.Range("H11").Value = 0
$wsPaper.Range("B12").Value = 9.376521199999999
$wsPaper.Range("C12").Value = 16.8232429
$wsPaper.Range("D12").Value = $null
$wsPaper.Range("E12").Value = 2.2839657
$wsPaper.Range("F12").Value = 3.6041693
$wsPaper.Range("G12").Value = 10.9351079
$wsPaper.Range("H12").Value = 0
$wsPaper.Range("B13").Value = -0.0159034
$wsPaper.Range("C13").Value = -0.009961899999999999
$wsPaper.Range("D13").Value = $null
$wsPaper.Range("E13").Value = -0.0013525
$wsPaper.Range("F13").Value = -0.0021342
$wsPaper.Range("G13").Value = -0.0064753
$wsPaper.Range("H13").Value = -0
$wsPaper.Range("B14").Value = 0.0361373
$wsPaper.Range("C14").Value = 0.0316076
$wsPaper.Range("D14").Value = $null
$wsPaper.Range("E14").Value = 0.0042911
$wsPaper.Range("F14").Value = 0.0067715
$wsPaper.Range("G14").Value = 0.0205449
$wsPaper.Range("H14").Value = 0
$wsPaper.Range("B15").Value = 0.6626444
$wsPaper.Range("C15").Value = 1.5051356
$wsPaper.Range("D15").Value = $null
$wsPaper.Range("E15").Value = 0.204341
$wsPaper.Range("F15").Value = 0.3224565
$wsPaper.Range("G15").Value = 0.9783381
$wsPaper.Range("H15").Value = 0
$wsPaper.Range("B16").Value = 3.2045807
$wsPaper.Range("C16").Value = 5.5627593
$wsPaper.Range("D16").Value = $null
$wsPaper.Range("E16").Value = 0.7552141999999999
$wsPaper.Range("F16").Value = 1.1917516
$wsPaper.Range("G16").Value = 3.6157936
$wsPaper.Range("H16").Value = 0
$wsPaper.Range("B17").Value = 5.3627504
$wsPaper.Range("D17").Value = $null
$wsPaper.Range("E17").Value = 2.2906761
$wsPaper.Range("F17").Value = 3.6147585
$wsPaper.Range("G17").Value = 10.9672357
$wsPaper.Range("H17").Value = 0
$wsPaper.Range("B18").Value = 2.8913937
$wsPaper.Range("C18").Value = 6.4776482
$wsPaper.Range("D18").Value = $null
$wsPaper.Range("E18").Value = 0.8794218
$wsPaper.Range("F18").Value = 1.3877551
$wsPaper.Range("G18").Value = 4.2104713
$wsPaper.Range("H18").Value = 0
$wsPaper.Range("B19").Value = 2.8393539
$wsPaper.Range("C19").Value = 9.8557275
$wsPaper.Range("D19").Value = $null
$wsPaper.Range("E19").Value = 1.3380383
$wsPaper.Range("F19").Value = 2.1114663
$wsPaper.Range("G19").Value = 6.4062229
$wsPaper.Range("H19").Value = 0
$wsPaper.Range("B20").Value = 0.6004588
$wsPaper.Range("D20").Value = $null
$wsPaper.Range("E20").Value = 0.1288816
$wsPaper.Range("F20").Value = 0.2033792
$wsPaper.Range("G20").Value = 0.6170557
$wsPaper.Range("H20").Value = 0
$wsPaper.Range("B21").Value = 0.6850908999999999
$wsPaper.Range("C21").Value = 1.6879142
$wsPaper.Range("D21").Value = $null
$wsPaper.Range("E21").Value = 0.2291555
$wsPaper.Range("F21").Value = 0.3616145
$wsPaper.Range("G21").Value = 1.0971442
$wsPaper.Range("H21").Value = 0
$wsPaper.Range("B22").Value = 1.1839562
$wsPaper.Range("C22").Value = 4.8129031
$wsPaper.Range("D22").Value = $null
$wsPaper.Range("E22").Value = 0.6534118
$wsPaper.Range("F22").Value = 1.0311043
$wsPaper.Range("G22").Value = 3.128387
$wsPaper.Range("H22").Value = 0
$wsPaper.Range("B23").Value = 24.1617104
$wsPaper.Range("C23").Value = 68.4785462
$wsPaper.Range("D23").Value = $null
$wsPaper.Range("E23").Value = 9.296819299999999
$wsPaper.Range("F23").Value = 14.6706719
$wsPaper.Range("G23").Value = 44.511055
$wsPaper.Range("H23").Value = 0
$wsPaper.Range("B24").Value = 23.4636222
$wsPaper.Range("C24").Value = 54.6406485
$wsPaper.Range("D24").Value = $null
$wsPaper.Range("E24").Value = 7.4181516
$wsPaper.Range("F24").Value = 11.7060754
$wsPaper.Range("G24").Value = 35.5164215
$wsPaper.Range("H24").Value = 0
$wsPaper.Range("B25").Value = 19.694248
$wsPaper.Range("D25").Value = $null
$wsPaper.Range("E25").Value = 2.2361119
$wsPaper.Range("F25").Value = 3.5286546
$wsPaper.Range("G25").Value = 10.7059949
$wsPaper.Range("H25").Value = 0
$wsPaper.Range("B26").Value = 6.4551637
$wsPaper.Range("C26").Value = 3.6505394
$wsPaper.Range("D26").Value = $null
$wsPaper.Range("E26").Value = 0.4956064
$wsPaper.Range("F26").Value = 0.7820824
$wsPaper.Range("G26").Value = 2.3728506
$wsPaper.Range("H26").Value = 0
$wsPaper.Range("B29").Value = 0.0050995
$wsPaper.Range("C29").Value = 0.0199243
$wsPaper.Range("D29").Value = $null
$wsPaper.Range("E29").Value = 0.002705
$wsPaper.Range("F29").Value = 0.0042685
$wsPaper.Range("G29").Value = 0.0129508
$wsPaper.Range("H29").Value = 0
$wsPaper.Range("B31").Value = 0.3399668
$wsPaper.Range("C31").Value = 0.6164207
$wsPaper.Range("D31").Value = $null
$wsPaper.Range("E31").Value = 0.08368680000000001
$wsPaper.Range("F31").Value = 0.1320604
$wsPaper.Range("G31").Value = 0.4006734
$wsPaper.Range("H31").Value = 0
$wsPaper.Range("B34").Value = 0.08228630000000001
$wsPaper.Range("C34").Value = -0.2111502
$wsPaper.Range("D34").Value = $null
$wsPaper.Range("E34").Value = -0.0286663
$wsPaper.Range("F34").Value = -0.0452363
$wsPaper.Range("G34").Value = -0.1372476
$wsPaper.Range("H34").Value = -0
$wsPaper.Range("B35").Value = 0.4280643
$wsPaper.Range("C35").Value = 0.2651854
$wsPaper.Range("D35").Value = $null
$wsPaper.Range("E35").Value = 0.0360022
$wsPaper.Range("F35").Value = 0.0568127
$wsPaper.Range("G35").Value = 0.1723705
$wsPaper.Range("H35").Value = 0

# forecasted sheet updates
$wsForecasted.Range("B2").Value = 12.3839154
$wsForecasted.Range("C2").Value = 66.53384579999999
$wsForecasted.Range("D2").Value = $null
$wsForecasted.Range("E2").Value = 3.4987237
$wsForecasted.Range("F2").Value = 6.1691571
$wsForecasted.Range("G2").Value = 12.4509913
$wsForecasted.Range("H2").Value = 44.4149737
$wsForecasted.Range("B3").Value = 2.3020918
$wsForecasted.Range("C3").Value = 7.5462938
$wsForecasted.Range("D3").Value = $null
$wsForecasted.Range("E3").Value = 0.6616338
$wsForecasted.Range("F3").Value = 0.6381361
$wsForecasted.Range("G3").Value = 1.6722188
$wsForecasted.Range("H3").Value = 4.5743051
$wsForecasted.Range("B4").Value = 5.7043363
$wsForecasted.Range("C4").Value = 37.9989046
$wsForecasted.Range("D4").Value = $null
$wsForecasted.Range("E4").Value = 1.9465626
$wsForecasted.Range("F4").Value = 2.8801129
$wsForecasted.Range("G4").Value = 8.875712399999999
$wsForecasted.Range("H4").Value = 24.2965167
$wsForecasted.Range("B5").Value = 1.4856908
$wsForecasted.Range("D5").Value = $null
$wsForecasted.Range("E5").Value = 0.1836844
$wsForecasted.Range("F5").Value = 0.1686277
$wsForecasted.Range("G5").Value = 0.5649303
$wsForecasted.Range("H5").Value = 2.5733415
$wsForecasted.Range("B6").Value = 79.34370010000001
$wsForecasted.Range("C6").Value = 314.0469573
$wsForecasted.Range("D6").Value = $null
$wsForecasted.Range("E6").Value = 16.785604
$wsForecasted.Range("F6").Value = 26.2474493
$wsForecasted.Range("G6").Value = 68.9609253
$wsForecasted.Range("H6").Value = 202.0529786
$wsForecasted.Range("B7").Value = 0.8215785
$wsForecasted.Range("C7").Value = 3.0442464
$wsForecasted.Range("D7").Value = $null
$wsForecasted.Range("E7").Value = 0.1093115
$wsForecasted.Range("F7").Value = 0.0263104
$wsForecasted.Range("G7").Value = 0.2065458
$wsForecasted.Range("H7").Value = 2.7020787
$wsForecasted.Range("B8").Value = 5.9602057
$wsForecasted.Range("C8").Value = 13.8998462
$wsForecasted.Range("D8").Value = $null
$wsForecasted.Range("E8").Value = 0.5842299
$wsForecasted.Range("F8").Value = 0.2566977
$wsForecasted.Range("G8").Value = 1.1909938
$wsForecasted.Range("H8").Value = 11.8679249
$wsForecasted.Range("B9").Value = 28.9045907
$wsForecasted.Range("C9").Value = 85.0727284
$wsForecasted.Range("D9").Value = $null
$wsForecasted.Range("E9").Value = 4.5785642
$wsForecasted.Range("F9").Value = 5.521403
$wsForecasted.Range("G9").Value = 14.4699897
$wsForecasted.Range("H9").Value = 60.5027715
$wsForecasted.Range("B10").Value = 32.7960241
$wsForecasted.Range("C10").Value = 112.2357994
$wsForecasted.Range("D10").Value = $null
$wsForecasted.Range("E10").Value = 6.3167092
$wsForecasted.Range("F10").Value = 10.1052859
$wsForecasted.Range("G10").Value = 23.6476437
$wsForecasted.Range("H10").Value = 72.1661606
$wsForecasted.Range("B11").Value = 0.916753
$wsForecasted.Range("C11").Value = 4.258178
$wsForecasted.Range("D11").Value = $null
$wsForecasted.Range("E11").Value = 0.1972236
$wsForecasted.Range("F11").Value = 0.2025925
$wsForecasted.Range("G11").Value = 0.3658334
$wsForecasted.Range("H11").Value = 3.4925285
$wsForecasted.Range("B12").Value = 35.1435568
$wsForecasted.Range("C12").Value = 120.7647884
$wsForecasted.Range("D12").Value = $null
$wsForecasted.Range("E12").Value = 6.2162332
$wsForecasted.Range("F12").Value = 8.0236368
$wsForecasted.Range("G12").Value = 23.5464832
$wsForecasted.Range("H12").Value = 82.9784353
$wsForecasted.Range("B13").Value = 0.4943123
$wsForecasted.Range("C13").Value = 1.3472278
$wsForecasted.Range("D13").Value = $null
$wsForecasted.Range("E13").Value = 0.0407459
$wsForecasted.Range("F13").Value = 0.0270124
$wsForecasted.Range("G13").Value = 0.1678777
$wsForecasted.Range("H13").Value = 1.1115918
$wsForecasted.Range("B14").Value = 2.5516678
$wsForecasted.Range("C14").Value = 3.6622113
$wsForecasted.Range("D14").Value = $null
$wsForecasted.Range("E14").Value = 0.1172298
$wsForecasted.Range("F14").Value = 0.0956694
$wsForecasted.Range("G14").Value = 0.462041
$wsForecasted.Range("H14").Value = 2.9872711
$wsForecasted.Range("B15").Value = 3.1543616
$wsForecasted.Range("C15").Value = 16.959397
$wsForecasted.Range("D15").Value = $null
$wsForecasted.Range("E15").Value = 0.9196776
$wsForecasted.Range("F15").Value = 1.6715441
$wsForecasted.Range("G15").Value = 2.8275782
$wsForecasted.Range("H15").Value = 11.540597
$wsForecasted.Range("B16").Value = 14.5121537
$wsForecasted.Range("C16").Value = 82.6566736
$wsForecasted.Range("D16").Value = $null
$wsForecasted.Range("E16").Value = 4.1863333
$wsForecasted.Range("F16").Value = 8.3200862
$wsForecasted.Range("G16").Value = 14.7179459
$wsForecasted.Range("H16").Value = 55.4323082
$wsForecasted.Range("B17").Value = 10.5653163
$wsForecasted.Range("D17").Value = $null
$wsForecasted.Range("E17").Value = 3.7679905
$wsForecasted.Range("F17").Value = 5.3651822
$wsForecasted.Range("G17").Value = 18.8776215
$wsForecasted.Range("H17").Value = 34.9208694
$wsForecasted.Range("B18").Value = 12.6592355
$wsForecasted.Range("C18").Value = 51.3836144
$wsForecasted.Range("D18").Value = $null
$wsForecasted.Range("E18").Value = 3.0348563
$wsForecasted.Range("F18").Value = 3.8833055
$wsForecasted.Range("G18").Value = 8.9216734
$wsForecasted.Range("H18").Value = 35.5437792
$wsForecasted.Range("B19").Value = 6.1567263
$wsForecasted.Range("C19").Value = 24.3063215
$wsForecasted.Range("D19").Value = $null
$wsForecasted.Range("E19").Value = 1.9463003
$wsForecasted.Range("F19").Value = 2.8074386
$wsForecasted.Range("G19").Value = 7.3586128
$wsForecasted.Range("H19").Value = 12.1939699
$wsForecasted.Range("B20").Value = 7.1939678
$wsForecasted.Range("D20").Value = $null
$wsForecasted.Range("E20").Value = 0.7542327
$wsForecasted.Range("F20").Value = 0.9076664
$wsForecasted.Range("G20").Value = 2.3632568
$wsForecasted.Range("H20").Value = 12.6685278
$wsForecasted.Range("B21").Value = 2.368295
$wsForecasted.Range("C21").Value = 3.3934732
$wsForecasted.Range("D21").Value = $null
$wsForecasted.Range("E21").Value = 0.304665
$wsForecasted.Range("F21").Value = 0.3911398
$wsForecasted.Range("G21").Value = 1.2217401
$wsForecasted.Range("H21").Value = 1.4759283
$wsForecasted.Range("B22").Value = 6.5495747
$wsForecasted.Range("C22").Value = 30.3335259
$wsForecasted.Range("D22").Value = $null
$wsForecasted.Range("E22").Value = 1.4723935
$wsForecasted.Range("F22").Value = 2.0043844
$wsForecasted.Range("G22").Value = 7.2989919
$wsForecasted.Range("H22").Value = 19.5577561
$wsForecasted.Range("B23").Value = 27.843328
$wsForecasted.Range("C23").Value = 90.8392281
$wsForecasted.Range("D23").Value = $null
$wsForecasted.Range("E23").Value = 10.1191722
$wsForecasted.Range("F23").Value = 15.6050521
$wsForecasted.Range("G23").Value = 48.1722093
$wsForecasted.Range("H23").Value = 16.9427944
$wsForecasted.Range("B24").Value = 29.52942
$wsForecasted.Range("C24").Value = 80.9968568
$wsForecasted.Range("D24").Value = $null
$wsForecasted.Range("E24").Value = 8.424507200000001
$wsForecasted.Range("F24").Value = 12.8341374
$wsForecasted.Range("G24").Value = 39.7612655
$wsForecasted.Range("H24").Value = 19.9769467
$wsForecasted.Range("B25").Value = 32.527615
$wsForecasted.Range("D25").Value = $null
$wsForecasted.Range("E25").Value = 4.1508933
$wsForecasted.Range("F25").Value = 5.4520615
$wsForecasted.Range("G25").Value = 19.8963846
$wsForecasted.Range("H25").Value = 47.5790176
$wsForecasted.Range("B26").Value = 30.8318984
$wsForecasted.Range("C26").Value = 10.0838673
$wsForecasted.Range("D26").Value = $null
$wsForecasted.Range("E26").Value = 1.0228523
$wsForecasted.Range("F26").Value = 1.4192778
$wsForecasted.Range("G26").Value = 2.9302969
$wsForecasted.Range("H26").Value = 4.7114403
$wsForecasted.Range("B29").Value = 0.2679165
$wsForecasted.Range("C29").Value = 0.9098912
$wsForecasted.Range("D29").Value = $null
$wsForecasted.Range("E29").Value = 0.0336294
$wsForecasted.Range("F29").Value = 0.0118051
$wsForecasted.Range("G29").Value = 0.0774854
$wsForecasted.Range("H29").Value = 0.7869712
$wsForecasted.Range("B31").Value = 1.3772429
$wsForecasted.Range("C31").Value = 9.319411799999999
$wsForecasted.Range("D31").Value = $null
$wsForecasted.Range("E31").Value = 0.4235017
$wsForecasted.Range("F31").Value = 0.6135906
$wsForecasted.Range("G31").Value = 1.5731574
$wsForecasted.Range("H31").Value = 6.7091621
$wsForecasted.Range("B34").Value = 0.3672481
$wsForecasted.Range("C34").Value = 2.190498
$wsForecasted.Range("D34").Value = $null
$wsForecasted.Range("E34").Value = 0.1033395
$wsForecasted.Range("F34").Value = 0.2368027
$wsForecasted.Range("G34").Value = 0.07792739999999999
$wsForecasted.Range("H34").Value = 1.7724284
$wsForecasted.Range("B35").Value = 0.5534422
$wsForecasted.Range("C35").Value = 1.0328745
$wsForecasted.Range("D35").Value = $null
$wsForecasted.Range("E35").Value = 0.0663185
$wsForecasted.Range("F35").Value = 0.0817596
$wsForecasted.Range("G35").Value = 0.2048333
$wsForecasted.Range("H35").Value = 0.6799631
